$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency symbol/price/volume data (GitHub Actions scrape, 2022-12-17 21:58 UTC).
# Values are stored as text in the sheet, so each cell's number format is forced to "@" (Text)
# before assignment - this stops Excel from re-interpreting numeric-looking strings (e.g. "236.83")
# as numbers - and the cell's original style is restored afterwards so formatting is untouched.
$updates = @{
    'D2' = '236.83'
    'D4' = '5.355'
    'D5' = '0.05562'
    'D6' = '3.369'
    'D7' = '6.460'
    'D8' = '0.7990'
    'D9' = '1.038'
    'D10' = '0.1399'
    'D11' = '0.07292'
    'D12' = '0.03207'
    'D14' = '0.09253'
    'D15' = '0.001669'
    'D16' = '3.252'
    'D17' = '0.04764'
    'D18' = '0.0005709'
    'E18' = '17OneONEWorstin24h'
    'D19' = '0.006266'
    'D20' = '0.005069'
    'D21' = '0.001048'
    'D22' = '0.0001499'
    'D23' = '0.0004181'
    'D24' = '3.942'
    'D27' = '0.1295'
    'D40' = '0.04136'
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D42' = '0.003499'
    'E42' = '41CEJICEJI'
    'B43' = 'BKEXToken'
    'C43' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D43' = '0.1037'
    'E43' = '42BKEXTokenBKK'
    'D44' = '0.008832'
    'D45' = '0.00005441'
    'D46' = '0.00000000750'
    'D47' = '0.6800'
    'D48' = '0.03200'
    'E48' = '47BOLOBOLO'
    'D49' = '0.00002099'
    'D50' = '0.01010'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = $origStyle
}

Write-Host "Applied $($updates.Count) cell updates"